$d = $word.ActiveDocument

# Locate the run "POST -D headers <en-dash>H " — this spans the unchanged
# "POST" run plus the run we need to edit/split (" -D headers \u2013H ").
# We only want to touch the single en-dash character inside the second run,
# turning it into a plain hyphen, while leaving "POST" (and everything else)
# completely untouched as separate runs.
$probe = $d.Content
$found = $probe.Find.Execute("POST -D headers –H ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate target text"
}

# The en-dash is the 3rd-from-last character of the matched range
# ("...H " are the final two characters, the dash sits right before them).
$dashStart = $probe.End - 3
$dashEnd = $probe.End - 2

# Step 1: force a clean run split around the single en-dash character using
# a character-formatting-only mutation (Bold on/off). Formatting-only edits
# split runs at the touched boundaries without merging same-formatted
# neighboring runs together (unlike a direct text/content edit, which would
# coalesce "POST" into the same run as the text we are changing).
$dashRange = $d.Range($dashStart, $dashEnd)
$dashRange.Bold = 1

# Step 2: with the en-dash now isolated in its own (bold) run, replace its
# text with a plain hyphen. Because this run's formatting still differs from
# its neighbors at this point, the text replacement cannot re-merge it with
# "POST" or the following "H " run.
$dashRange2 = $d.Range($dashStart, $dashEnd)
$dashRange2.Text = "-"

# Step 3: revert the temporary bold formatting back to normal. Pure
# formatting changes do not trigger the run-coalescing behavior, so the
# three-way split (" -D headers ", "-", "H ") survives intact.
$dashRange3 = $d.Range($dashStart, $dashEnd)
$dashRange3.Bold = 0
